$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24 (participant #23): Kamga Jaures
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "Kamga Jaures"
$ws.Range("C24").Value = 691348635
$ws.Range("D24").Value = "M"
$ws.Range("E24").Value = 500
$ws.Range("F24").Value = "cash"
$ws.Range("G24").Value = 200

# Row 25 (participant #24): Nsini Franc
$ws.Range("A25").Value = 24
$ws.Range("B25").Value = "Nsini Franc"
$ws.Range("C25").Value = 658253270
$ws.Range("D25").Value = "M"
$ws.Range("E25").Value = 500
$ws.Range("F25").Value = "cash"
$ws.Range("G25").Value = 0

# Extend the "Amount left to reimburse" shared formula down into the new rows
$ws.Range("H24:H25").Formula = "=E24-G24-300"

# Match the updated view/selection state from the edit
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
$ws.Range("H26").Select()
